$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Populate the "100 Samples" worksheet with the full 100-sample dataset
#    (previously empty). Build a 2D array and write it in one shot.
# ---------------------------------------------------------------------------
$ws100 = $wb.Worksheets.Item("100 Samples")

$arr = New-Object 'object[,]' 30,17
$arr[0,1] = 'Target.Maneuver'
$arr[0,2] = 'Target.Quantity'
$arr[0,3] = 'Launch.Range'
$arr[0,4] = 'Approach.Angle'
$arr[0,5] = 'Ship.Maneuver'
$arr[0,6] = 'day'
$arr[0,7] = 'wholeplot_error'
$arr[0,8] = 'subplot_error'
$arr[0,11] = 'Estimate'
$arr[0,12] = 'Std..Error'
$arr[0,13] = 't.value'
$arr[0,14] = 'p.z'
$arr[0,15] = 'Power'
$arr[0,16] = 'RunsPerSim'
$arr[1,0] = 2
$arr[1,1] = 'Radial'
$arr[1,2] = 3
$arr[1,3] = 2500
$arr[1,4] = 0
$arr[1,5] = 'No'
$arr[1,6] = 1
$arr[1,7] = -0.812182442683306
$arr[1,8] = 1.70307363179695
$arr[1,10] = '(Intercept)'
$arr[1,11] = 0.586445567579611
$arr[1,12] = 0.249818128065514
$arr[1,13] = 2.34749004053789
$arr[1,14] = 0.0189003784403063
$arr[1,15] = 0.57
$arr[1,16] = 100
$arr[2,0] = 4
$arr[2,1] = 'Radial'
$arr[2,2] = 9
$arr[2,3] = 2500
$arr[2,4] = 45
$arr[2,5] = 'No'
$arr[2,6] = 1
$arr[2,7] = -0.812182442683306
$arr[2,8] = 0.691186053997328
$arr[2,10] = 'design$Target.Maneuver'
$arr[2,11] = 0.579673427991478
$arr[2,12] = 0.388054730459349
$arr[2,13] = 1.49379296911367
$arr[2,14] = 0.135229739789927
$arr[2,15] = 0.91
$arr[2,16] = 100
$arr[3,0] = 5
$arr[3,1] = 'Radial'
$arr[3,2] = 9
$arr[3,3] = 5000
$arr[3,4] = 45
$arr[3,5] = 'No'
$arr[3,6] = 1
$arr[3,7] = -0.812182442683306
$arr[3,8] = 0.99688005885675
$arr[3,10] = 'design$Target.Quantity'
$arr[3,11] = 0.455922737718646
$arr[3,12] = 0.427605778664621
$arr[3,13] = 1.06622211501083
$arr[3,14] = 0.286323247391978
$arr[3,15] = 0.79
$arr[3,16] = 100
$arr[4,0] = 6
$arr[4,1] = 'Radial'
$arr[4,2] = 3
$arr[4,3] = 2500
$arr[4,4] = 45
$arr[4,5] = 'No'
$arr[4,6] = 1
$arr[4,7] = -0.812182442683306
$arr[4,8] = 1.69279581924356
$arr[4,10] = 'design$Launch.Range'
$arr[4,11] = -0.0204383539041419
$arr[4,12] = 0.528334298265064
$arr[4,13] = -0.0386845108698357
$arr[4,14] = 0.969141922701087
$arr[4,15] = 0.74
$arr[4,16] = 100
$arr[5,0] = 7
$arr[5,1] = 'Radial'
$arr[5,2] = 9
$arr[5,3] = 7500
$arr[5,4] = 0
$arr[5,5] = 'No'
$arr[5,6] = 2
$arr[5,7] = -0.812182442683306
$arr[5,8] = 0.255249114524872
$arr[5,10] = 'design$Approach.Angle'
$arr[5,11] = -0.0365043211351332
$arr[5,12] = 0.468398174954159
$arr[5,13] = -0.0779343795238012
$arr[5,14] = 0.937880251662626
$arr[5,15] = 0.71
$arr[5,16] = 100
$arr[6,0] = 8
$arr[6,1] = 'Weave'
$arr[6,2] = 9
$arr[6,3] = 5000
$arr[6,4] = 45
$arr[6,5] = 'No'
$arr[6,6] = 1
$arr[6,7] = -0.812182442683306
$arr[6,8] = -2.88618258475824
$arr[6,10] = 'design$Ship.Maneuver'
$arr[6,11] = -0.365612060220531
$arr[6,12] = 0.612981620820512
$arr[6,13] = -0.596448650011949
$arr[6,14] = 0.550875548282684
$arr[6,15] = 0.59
$arr[6,16] = 100
$arr[7,0] = 9
$arr[7,1] = 'Weave'
$arr[7,2] = 9
$arr[7,3] = 2500
$arr[7,4] = 45
$arr[7,5] = 'No'
$arr[7,6] = 1
$arr[7,7] = -0.812182442683306
$arr[7,8] = 0.496489998521541
$arr[7,10] = 'design$Target.Maneuver:design$Target.Quantity'
$arr[7,11] = 0.643956689265206
$arr[7,12] = 0.748913980275151
$arr[7,13] = 0.859854010241091
$arr[7,14] = 0.389869522438307
$arr[7,15] = 0.44
$arr[7,16] = 100
$arr[8,0] = 10
$arr[8,1] = 'Weave'
$arr[8,2] = 3
$arr[8,3] = 2500
$arr[8,4] = 45
$arr[8,5] = 'No'
$arr[8,6] = 1
$arr[8,7] = -0.812182442683306
$arr[8,8] = -0.670264584540183
$arr[8,10] = 'design$Target.Maneuver:design$Launch.Range'
$arr[8,11] = 1.21744050309255
$arr[8,12] = 1.11408616640454
$arr[8,13] = 1.09277050537443
$arr[8,14] = 0.274494574633541
$arr[8,15] = 0.33
$arr[8,16] = 100
$arr[9,0] = 11
$arr[9,1] = 'Radial'
$arr[9,2] = 3
$arr[9,3] = 5000
$arr[9,4] = 0
$arr[9,5] = 'No'
$arr[9,6] = 1
$arr[9,7] = -2.67153526106257
$arr[9,8] = 0.60638227013168
$arr[9,10] = 'design$Target.Maneuver:design$Approach.Angle'
$arr[9,11] = -0.208276046213262
$arr[9,12] = 0.848801547424425
$arr[9,13] = -0.245376609933439
$arr[9,14] = 0.806164838071003
$arr[9,15] = 0.39
$arr[9,16] = 100
$arr[10,0] = 13
$arr[10,1] = 'Weave'
$arr[10,2] = 9
$arr[10,3] = 2500
$arr[10,4] = 0
$arr[10,5] = 'No'
$arr[10,6] = 1
$arr[10,7] = -2.67153526106257
$arr[10,8] = -1.95298373104632
$arr[10,10] = 'design$Target.Quantity:design$Launch.Range'
$arr[10,11] = 0.415621806511803
$arr[10,12] = 1.20331855935541
$arr[10,13] = 0.345396323592351
$arr[10,14] = 0.729796440664115
$arr[10,15] = 0.32
$arr[10,16] = 100
$arr[11,0] = 14
$arr[11,1] = 'Weave'
$arr[11,2] = 3
$arr[11,3] = 5000
$arr[11,4] = 0
$arr[11,5] = 'No'
$arr[11,6] = 1
$arr[11,7] = -2.67153526106257
$arr[11,8] = 1.48326413654553
$arr[11,10] = 'design$Target.Quantity:design$Ship.Maneuver'
$arr[11,11] = 0.117107821925629
$arr[11,12] = 0.900552606805264
$arr[11,13] = 0.130039956622936
$arr[11,14] = 0.89653481425799
$arr[11,15] = 0.45
$arr[11,16] = 100
$arr[12,0] = 17
$arr[12,1] = 'Radial'
$arr[12,2] = 9
$arr[12,3] = 7500
$arr[12,4] = 0
$arr[12,5] = 'No'
$arr[12,6] = 1
$arr[12,7] = -2.67153526106257
$arr[12,8] = 0.828455715182223
$arr[13,0] = 18
$arr[13,1] = 'Weave'
$arr[13,2] = 3
$arr[13,3] = 2500
$arr[13,4] = 0
$arr[13,5] = 'No'
$arr[13,6] = 2
$arr[13,7] = -2.67153526106257
$arr[13,8] = 1.17942876910367
$arr[14,0] = 20
$arr[14,1] = 'Radial'
$arr[14,2] = 3
$arr[14,3] = 7500
$arr[14,4] = 45
$arr[14,5] = 'No'
$arr[14,6] = 1
$arr[14,7] = -2.67153526106257
$arr[14,8] = 0.810319243202051
$arr[15,0] = 21
$arr[15,1] = 'Weave'
$arr[15,2] = 9
$arr[15,3] = 7500
$arr[15,4] = 45
$arr[15,5] = 'No'
$arr[15,6] = 1
$arr[15,7] = -2.67153526106257
$arr[15,8] = -0.199347221982706
$arr[16,0] = 23
$arr[16,1] = 'Radial'
$arr[16,2] = 9
$arr[16,3] = 5000
$arr[16,4] = 45
$arr[16,5] = 'No'
$arr[16,6] = 1
$arr[16,7] = -2.67153526106257
$arr[16,8] = 0.860336346706679
$arr[17,0] = 24
$arr[17,1] = 'Weave'
$arr[17,2] = 3
$arr[17,3] = 2500
$arr[17,4] = 45
$arr[17,5] = 'No'
$arr[17,6] = 1
$arr[17,7] = -2.67153526106257
$arr[17,8] = -1.42237393094252
$arr[18,0] = 25
$arr[18,1] = 'Weave'
$arr[18,2] = 9
$arr[18,3] = 5000
$arr[18,4] = 0
$arr[18,5] = 'No'
$arr[18,6] = 1
$arr[18,7] = -2.67153526106257
$arr[18,8] = 0.373062278243411
$arr[19,0] = 26
$arr[19,1] = 'Weave'
$arr[19,2] = 3
$arr[19,3] = 5000
$arr[19,4] = 45
$arr[19,5] = 'No'
$arr[19,6] = 1
$arr[19,7] = -2.67153526106257
$arr[19,8] = 0.639543363329221
$arr[20,0] = 27
$arr[20,1] = 'Radial'
$arr[20,2] = 3
$arr[20,3] = 5000
$arr[20,4] = 45
$arr[20,5] = 'No'
$arr[20,6] = 1
$arr[20,7] = -2.67153526106257
$arr[20,8] = -0.259507559250752
$arr[21,0] = 1
$arr[21,1] = 'Weave'
$arr[21,2] = 3
$arr[21,3] = 7500
$arr[21,4] = 0
$arr[21,5] = 'Yes'
$arr[21,6] = 2
$arr[21,7] = -2.67153526106257
$arr[21,8] = -0.145099634574724
$arr[22,0] = 3
$arr[22,1] = 'Radial'
$arr[22,2] = 3
$arr[22,3] = 5000
$arr[22,4] = 0
$arr[22,5] = 'Yes'
$arr[22,6] = 2
$arr[22,7] = -2.67153526106257
$arr[22,8] = 1.46786727541917
$arr[23,0] = 12
$arr[23,1] = 'Radial'
$arr[23,2] = 3
$arr[23,3] = 2500
$arr[23,4] = 0
$arr[23,5] = 'Yes'
$arr[23,6] = 2
$arr[23,7] = -2.67153526106257
$arr[23,8] = -0.529016892051686
$arr[24,0] = 15
$arr[24,1] = 'Weave'
$arr[24,2] = 3
$arr[24,3] = 7500
$arr[24,4] = 0
$arr[24,5] = 'Yes'
$arr[24,6] = 2
$arr[24,7] = -2.67153526106257
$arr[24,8] = -0.087945411630587
$arr[25,0] = 16
$arr[25,1] = 'Radial'
$arr[25,2] = 9
$arr[25,3] = 5000
$arr[25,4] = 0
$arr[25,5] = 'Yes'
$arr[25,6] = 2
$arr[25,7] = -2.67153526106257
$arr[25,8] = -1.92721299046048
$arr[26,0] = 19
$arr[26,1] = 'Weave'
$arr[26,2] = 9
$arr[26,3] = 5000
$arr[26,4] = 0
$arr[26,5] = 'Yes'
$arr[26,6] = 2
$arr[26,7] = -2.67153526106257
$arr[26,8] = -0.933796150627326
$arr[27,0] = 22
$arr[27,1] = 'Radial'
$arr[27,2] = 9
$arr[27,3] = 2500
$arr[27,4] = 0
$arr[27,5] = 'Yes'
$arr[27,6] = 2
$arr[27,7] = -2.67153526106257
$arr[27,8] = 1.30045869196953
$arr[28,0] = 28
$arr[28,1] = 'Weave'
$arr[28,2] = 9
$arr[28,3] = 2500
$arr[28,4] = 0
$arr[28,5] = 'Yes'
$arr[28,6] = 2
$arr[28,7] = -2.67153526106257
$arr[28,8] = -0.628668226049944
$arr[29,0] = 29
$arr[29,1] = 'Weave'
$arr[29,2] = 3
$arr[29,3] = 7500
$arr[29,4] = 0
$arr[29,5] = 'Yes'
$arr[29,6] = 2
$arr[29,7] = -2.67153526106257
$arr[29,8] = -2.58706798978451

$ws100.Range("A1:Q30").Value2 = $arr

# ---------------------------------------------------------------------------
# 2. View-state updates that accompanied the data population:
#    - "100 Samples" becomes the active/selected tab (was "10 Samples").
#    - Each of the three "Samples" sheets gets its remembered cell selection
#      updated.
#    - Column F on "1 Sample" gets an explicit autofit width (it now holds
#      the same "Yes"/"No" Ship.Maneuver text as the other samples sheets).
#    - "LF21 FIREX DOE" scrolls so row 2 is the top visible row.
# ---------------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item("1 Sample")
$ws1.Columns.Item(6).AutoFit()
$ws1.Activate()
$ws1.Range("K14").Select()

$ws10 = $wb.Worksheets.Item("10 Samples")
$ws10.Activate()
$ws10.Range("K16").Select()

$ws100.Activate()
$ws100.Range("K14").Select()

$wsDoe = $wb.Worksheets.Item("LF21 FIREX DOE")
$wsDoe.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$wsDoe.Range("L14").Select()

# "100 Samples" is the tab that should end up active/selected.
$ws100.Activate()
